$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 291, shifting existing rows (291..384) down to (292..385)
$ws.Rows.Item(291).Insert()

# Populate the newly inserted row 291 with the new data record
$ws.Cells.Item(291, 1).Value = 5
$ws.Cells.Item(291, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(291, 3).Value = "Maule"
$ws.Cells.Item(291, 4).Value = 45215
$ws.Cells.Item(291, 5).Value = 7
$ws.Cells.Item(291, 6).Value = 100112017
$ws.Cells.Item(291, 7).Value = "Apio"
$ws.Cells.Item(291, 8).Value = "Americana (o)"
$ws.Cells.Item(291, 9).Value = "Primera"
$ws.Cells.Item(291, 10).Value = 500
$ws.Cells.Item(291, 11).Value = 7500
$ws.Cells.Item(291, 12).Value = 7500
$ws.Cells.Item(291, 13).Value = 7500
$ws.Cells.Item(291, 14).Value = "`$/docena de matas"
$ws.Cells.Item(291, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(291, 16).Value = 1250
$ws.Cells.Item(291, 17).Value = 6
$ws.Cells.Item(291, 18).Value = "Hortaliza"
